$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing C column values (remove unwanted columns from BY YEAR tables)
$ws.Range("C2").Value = 21
$ws.Range("C3").Value = 4
$ws.Range("C4").Value = 71
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 56
$ws.Range("C7").Value = 12
$ws.Range("C8").Value = 38

# Add a new row 9 to the 2021 group, copying the formatting from row 8
# so the new cells match the existing table styling (border/alignment/font)
$ws.Range("A8:B8").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Merge the year cell across the (now two-row) 2021 group
$ws.Range("A8:A9").Merge()

# Merging can re-derive border formatting on the merged range, so
# re-apply the original formatting afterward to keep it consistent
# with the other year groups
$ws.Range("A6:A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row's values
$ws.Range("A9").Value = $null
$ws.Range("B9").Value = $true
$ws.Range("C9").Value = 4
